$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A46").Value = "2023-12-06 21:18:03"
$ws.Range("B46").Value = 0.0004
